$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-9 from 2023-10-13 (45212)
# to 2023-10-22 (45221), preserving existing cell formatting.
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
